$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: was "Mon 23 Oct" data, now becomes old "Tue 24 Oct" data (as numbers)
$ws.Range("A2").Value = "Tue 24 Oct"
$ws.Range("B2").Value = 3113
$ws.Range("C2").Value = 10983
$ws.Range("D2").Value = 13
$ws.Range("E2").Value = 8.84

# Row 3: was "Tue 24 Oct" data, now becomes old "Wed 25 Oct" data (as numbers)
$ws.Range("A3").Value = "Wed 25 Oct"
$ws.Range("B3").Value = 2863
$ws.Range("C3").Value = 13708
$ws.Range("D3").Value = 26
$ws.Range("E3").Value = 10.88

# Row 4: was "Wed 25 Oct" data, now becomes a brand new "Thu 26 Oct" row
$ws.Range("A4").Value = "Thu 26 Oct"
$ws.Range("B4").Value = 2858
$ws.Range("C4").Value = 10185
$ws.Range("D4").Value = 19
$ws.Range("E4").Value = 8.28

# Row 5: "Today" row stays but with new numeric values
$ws.Range("A5").Value = "Today"
$ws.Range("B5").Value = 569
$ws.Range("C5").Value = 68
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0.05
